$wb = $excel.ActiveWorkbook

# Add the new worksheet after the current last sheet, matching the
# position/sheetId/rId ordering used by the target workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "two_by_two_PriceinDem"

# Row/column data for the new sheet (A1:I28). Text cells are prefixed with
# a leading apostrophe so the engine stores them as shared-string text
# (matching the quotePrefix cell style used throughout this workbook);
# $null marks a cell that should stay empty.
$rows = @(
    @($null, "'benchmark", "'RA=157", "'eRA=.5", "'pr_Ud=2", "'prU2,eRA.6", "'prU.5,eRA.6", "'Itax=0.1", "'Otax=0.1"),
    @("'X", 1, 1.0363877164248041, 1.0365519507521972, 1.0595609099982088, 1.0633532225145439, 0.99833470833445848, 0.98248710709607423, 0.94048657015713999),
    @("'Y", 1, 1.0432700717660239, 1.0430251457731119, 1.0087067867541124, 1.0030495908339756, 1.1000067470042523, 1.122931737696433, 1.1855074833626025),
    @("'U", 1, 1.0388246300487494, 1.0389899529717939, 1.0621536850625561, 1.0659719353994179, 1.0005254867593421, 0.99843756359281499, 0.99193526198210002),
    @("'PX", 1, 1.0023513532505028, 1.0023520310092886, 1.0024470278582813, 1.0024626933265708, 1.0021944327932055, 1.0162347743614695, 1.054704334391892),
    @("'PY", 1, 0.99573893485424703, 0.99573771069529127, 0.99556615181885155, 0.99553786527941868, 0.99602242514969408, 0.97114380037066672, 0.90770374205596749),
    @("'PU", 1, 1, 1, 1, 1, 1, 1, 1),
    @("'PL", 1, 0.94438602708152408, 0.94437046939986447, 0.94219250114219932, 0.94183384552627059, 0.9479952509302686, 0.89903690333116493, 0.84467427236056358),
    @("'PK", 1, 1.0388246300416317, 1.0388360220635884, 1.0404339284158413, 1.0406976527852307, 1.0361901997164102, 1.0329715504842176, 0.96149717694622527),
    @("'SX", 80, 80, 80, 80, 80, 80, 80, 80),
    @("'SY", 54, 54, 54, 54, 54, 54, 54, 54.000000000000007),
    @("'SU", 124, 124, 123.99999999999999, 124, 124, 124, 124, 124),
    @("'DXL", 30, 31.841365430850797, 31.841911518828791, 31.918541910798012, 31.931195659031214, 31.715172575173078, 30.827982414866039, 30.648734246825121),
    @("'DXK", 50, 48.244493065345502, 48.243996628540472, 48.174468386692752, 48.16301308279445, 48.359578823666638, 49.189872358276297, 49.36228226730212),
    @("'DYL", 24, 25.305048730086334, 25.305434497175867, 25.359560402662385, 25.368496662333634, 25.215883919305377, 25.924910448641455, 25.790876462309562),
    @("'DYK", 30, 28.755737186305545, 28.755386493763073, 28.706277005056187, 28.698187104056096, 28.837054010514596, 28.204372131496182, 28.321572766511753),
    @("'DUX", 80, 79.812333011330708, 79.812279044731554, 79.804715637612503, 79.803468530612463, 79.824829775826657, 78.721966634399251, 75.850641162032716),
    @("'DUY", 44, 44.188289178863421, 44.18834350389114, 44.195958168740567, 44.197213922797879, 44.175712201847752, 45.307399360636374, 48.473965635901294),
    @("'RA", 134, 139.20250041197306, 139.202487647439, 139.20094884111393, 139.20074264707893, 139.20613388257078, 138.76352763367868, 137.46365821827374),
    @("'DU", 124, 128.81425411257209, 128.83475414329422, 131.70705694774205, 132.18051998952777, 124.06516035720092, 123.80625788538022, 122.99997248577242),
    @("'DY", 10, 10.432700716801408, 10.412112941775911, 7.5272666509210593, 7.0516882404877128, 15.201438384376063, 15.401704405248296, 15.934368299221509),
    @("'CWI", 1, 1.0391557242667655, 1.0391553971951093, 1.0454956618531788, 1.046614299839608, 1.0534806670456223, 1.0530112713727526, 1.0506793942727906),
    @("'PX/PX", 1, 1, 1, 1, 1, 1, 1, 1),
    @("'PY/PX", 1, 0.99340309326184617, 0.99340120026759737, 0.99313592055419553, 0.99309218378573993, 0.9938415067560199, 0.95562937312479324, 0.86062388525151912),
    @("'PU/PX", 1, 0.99765416264179452, 0.997653488059559, 0.99755894547015678, 0.9975433566326557, 0.99781037219784841, 0.98402458292999251, 0.94813301452540932),
    @("'PL/PX", 1, 0.94217065105862896, 0.94215449281721775, 0.93989255786930181, 0.93952009565651817, 0.94591949417252408, 0.88467441383912149, 0.80086356414527782),
    @("'PK/PX", 1, 1.0363877164158561, 1.0363983809336559, 1.0378941724618791, 1.0381410297991052, 1.0339213288467941, 1.0164693991437799, 0.91162721683569548),
    @("'RA/PX", 134, 138.87595398615105, 138.87584734803519, 138.86115173438685, 138.85877606592558, 138.90132426159147, 136.54672240562513, 130.33383265418243)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $val = $row[$c]
        if ($val -ne $null) {
            $ws.Cells.Item($r + 1, $c + 1).Value = $val
        }
    }
}

# Make the new sheet the active tab (matches activeTab moving to this sheet).
$ws.Activate()
